$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data rows (2-5) with new credentials / candidate IDs
$ws.Range("A2").Value = "lIJgc882"
$ws.Range("B2").Value = 231006216
$ws.Range("C2").Value = "cycxjtt77"
$ws.Range("D2").Value = "a&`$A3p7E"

$ws.Range("A3").Value = "kQWks262"
$ws.Range("B3").Value = 231006215
$ws.Range("C3").Value = "gkkkjqj87"
$ws.Range("D3").Value = "sXF%3g2&"

$ws.Range("A4").Value = "bxKYK840"
$ws.Range("B4").Value = 231006214
$ws.Range("C4").Value = "cnhskzy54"
$ws.Range("D4").Value = "Qwy2`$9%M"

$ws.Range("A5").Value = "XPCuN181"
$ws.Range("B5").Value = 231006213
$ws.Range("C5").Value = "zwploqj61"
$ws.Range("D5").Value = "zk!64M#U"

# New row 6 with a brand new candidate record
$ws.Range("A6").Value = "uXYmz219"
$ws.Range("B6").Value = 231006212
$ws.Range("C6").Value = "wootqrz97"
$ws.Range("D6").Value = "SNe6&7%q"
$ws.Range("E6").Value = "MR"
$ws.Range("F6").Value = "lltIEWTQ"
$ws.Range("G6").Value = "dPlx"
$ws.Range("H6").Value = "Candidate"

# Match the formatting used by the other data rows
$ws.Range("A6:H6").Style = $ws.Range("A5:H5").Style
